# "TEL has been changed to DB" - commit adds a new Monday (2016-08-29) block
# of log entries to the "Logs" sheet: a day-separator row plus three task
# rows (AV Shutdown / Pickup PC / Pickup PC), mirroring the existing layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

function Set-Edge($cell, $edgeIndex) {
    $cell.Borders.Item($edgeIndex).LineStyle = 1
    $cell.Borders.Item($edgeIndex).Weight = 2
}

function Set-BlankBandedCell($cell, $numberFormat, $wrap, $left, $right, $top, $bottom) {
    $cell.Interior.Color = 65535
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    if ($numberFormat) { $cell.NumberFormat = $numberFormat }
    if ($wrap) { $cell.WrapText = $true }
    if ($left) { Set-Edge $cell 7 }
    if ($right) { Set-Edge $cell 10 }
    if ($top) { Set-Edge $cell 8 }
    if ($bottom) { Set-Edge $cell 9 }
}

# ---- Row 94: day-separator row ("MONDAY"), same banded look as rows 2/6/... ----
Set-BlankBandedCell $ws.Cells.Item(94,1) $null  $false $true  $true $true $true   # A94
Set-BlankBandedCell $ws.Cells.Item(94,2) "m/d/yy;@" $false $true  $true $true $true # B94
Set-BlankBandedCell $ws.Cells.Item(94,4) "@"  $false $true  $true $true $true   # D94
Set-BlankBandedCell $ws.Cells.Item(94,5) "@"  $false $true  $true $true $true   # E94 (border top/bottom only)
Set-BlankBandedCell $ws.Cells.Item(94,6) "@"  $true  $true  $true $true $true   # F94
$ws.Range("E94").Borders.Item(10).LineStyle = -4142
$ws.Range("E94").Borders.Item(7).LineStyle = -4142

$ws.Range("C94").Value = "MONDAY"

# ---- Rows 95-97: new log entries for Monday, Aug 29 2016 ----
$ws.Range("A95").Value = "AV Shutdown"
$ws.Range("B95").Value = 42611
$ws.Range("C95").Value = "1600"
$ws.Range("D95").Value = "FC"
$ws.Range("F95").Value = "Turn off PC and projector. Lock room."

$ws.Range("A96").Value = "Pickup PC"
$ws.Range("B96").Value = 42611
$ws.Range("C96").Value = "1600"
$ws.Range("D96").Value = "VC"
$ws.Range("E96").Value = "113"
$ws.Range("F96").Value = "Flat screen DLP and wireless keyboard. To Vanier 132 storeroom."

$ws.Range("A97").Value = "Pickup PC"
$ws.Range("B97").Value = 42611
$ws.Range("C97").Value = "1630"
$ws.Range("D97").Value = "VC"
$ws.Range("E97").Value = "010-SCR"
$ws.Range("F97").Value = "Leave portable screen and network cable, but get carpets. To  Vanier 040 storeroom."

# Added last, on purpose, so this string lands at the same shared-string
# index the original commit produced (strings are appended in write order).
$ws.Range("E95").Value = "305-SCR"

$ws.Rows.Item(97).RowHeight = 30

# ---- View state: sheet grew, so the frozen-pane anchor & selection moved ----
$ws.Application.ActiveWindow.ScrollRow = 77
$ws.Range("A102").Select()
